# Populate the (previously all-zero) team transition probability matrix
# with the newly simulated game outcome frequencies, per state row
# (rows 2-19, columns B-S). Cells not listed below remain 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2105263157894737
$ws.Cells.Item(2, 3).Value = 0.5263157894736842
$ws.Cells.Item(2, 10).Value = 0.05263157894736842
$ws.Cells.Item(2, 16).Value = 0.1052631578947368
$ws.Cells.Item(2, 19).Value = 0.1052631578947368
$ws.Cells.Item(3, 10).Value = 0.2
$ws.Cells.Item(3, 16).Value = 0.3
$ws.Cells.Item(3, 19).Value = 0.5
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(6, 6).Value = 0.07142857142857142
$ws.Cells.Item(6, 10).Value = 0.2142857142857143
$ws.Cells.Item(6, 17).Value = 0.1428571428571428
$ws.Cells.Item(6, 18).Value = 0.1428571428571428
$ws.Cells.Item(6, 19).Value = 0.4285714285714285
$ws.Cells.Item(7, 2).Value = 0.125
$ws.Cells.Item(7, 6).Value = 0.0625
$ws.Cells.Item(7, 10).Value = 0.125
$ws.Cells.Item(7, 17).Value = 0.1875
$ws.Cells.Item(7, 18).Value = 0.0625
$ws.Cells.Item(7, 19).Value = 0.4375
$ws.Cells.Item(8, 2).Value = 0.05
$ws.Cells.Item(8, 6).Value = 0.05
$ws.Cells.Item(8, 10).Value = 0.15
$ws.Cells.Item(8, 15).Value = 0.025
$ws.Cells.Item(8, 17).Value = 0.175
$ws.Cells.Item(8, 18).Value = 0.075
$ws.Cells.Item(8, 19).Value = 0.475
$ws.Cells.Item(9, 10).Value = 0.1333333333333333
$ws.Cells.Item(9, 15).Value = 0.06666666666666667
$ws.Cells.Item(9, 17).Value = 0.1333333333333333
$ws.Cells.Item(9, 19).Value = 0.6666666666666666
$ws.Cells.Item(10, 2).Value = 0.09322033898305085
$ws.Cells.Item(10, 4).Value = 0.03389830508474576
$ws.Cells.Item(10, 6).Value = 0.05084745762711865
$ws.Cells.Item(10, 10).Value = 0.2033898305084746
$ws.Cells.Item(10, 15).Value = 0.0423728813559322
$ws.Cells.Item(10, 17).Value = 0.1186440677966102
$ws.Cells.Item(10, 18).Value = 0.0423728813559322
$ws.Cells.Item(10, 19).Value = 0.4152542372881356
$ws.Cells.Item(11, 7).Value = 0.1818181818181818
$ws.Cells.Item(11, 11).Value = 0.2272727272727273
$ws.Cells.Item(11, 12).Value = 0.5454545454545454
$ws.Cells.Item(11, 19).Value = 0.04545454545454546
$ws.Cells.Item(12, 7).Value = 0.75
$ws.Cells.Item(12, 10).Value = 0.1666666666666667
$ws.Cells.Item(12, 19).Value = 0.08333333333333333
$ws.Cells.Item(13, 7).Value = 0.6
$ws.Cells.Item(13, 10).Value = 0.4
$ws.Cells.Item(15, 8).Value = 0.1666666666666667
$ws.Cells.Item(15, 9).Value = 0.05555555555555555
$ws.Cells.Item(15, 10).Value = 0.2777777777777778
$ws.Cells.Item(15, 13).Value = 0.05555555555555555
$ws.Cells.Item(15, 15).Value = 0.05555555555555555
$ws.Cells.Item(15, 19).Value = 0.3888888888888889
$ws.Cells.Item(16, 8).Value = 0.1111111111111111
$ws.Cells.Item(16, 10).Value = 0.4444444444444444
$ws.Cells.Item(16, 11).Value = 0.2222222222222222
$ws.Cells.Item(16, 13).Value = 0.1111111111111111
$ws.Cells.Item(16, 19).Value = 0.1111111111111111
$ws.Cells.Item(17, 6).Value = 0.03571428571428571
$ws.Cells.Item(17, 8).Value = 0.1428571428571428
$ws.Cells.Item(17, 9).Value = 0.1428571428571428
$ws.Cells.Item(17, 10).Value = 0.3214285714285715
$ws.Cells.Item(17, 11).Value = 0.1428571428571428
$ws.Cells.Item(17, 13).Value = 0.03571428571428571
$ws.Cells.Item(17, 15).Value = 0.07142857142857142
$ws.Cells.Item(17, 19).Value = 0.1071428571428571
$ws.Cells.Item(18, 8).Value = 0.09090909090909091
$ws.Cells.Item(18, 10).Value = 0.7272727272727273
$ws.Cells.Item(18, 19).Value = 0.1818181818181818
$ws.Cells.Item(19, 6).Value = 0.01587301587301587
$ws.Cells.Item(19, 8).Value = 0.246031746031746
$ws.Cells.Item(19, 9).Value = 0.07936507936507936
$ws.Cells.Item(19, 10).Value = 0.4047619047619048
$ws.Cells.Item(19, 11).Value = 0.0873015873015873
$ws.Cells.Item(19, 13).Value = 0.01587301587301587
$ws.Cells.Item(19, 15).Value = 0.04761904761904762
$ws.Cells.Item(19, 19).Value = 0.1031746031746032
